# Update the section description / title text boxes on slides 7-13.
# Each affected slide has shape index 8 ("TextBox 7") holding the longer
# descriptive sentence and shape index 9 ("TextBox 8") holding the short
# section label.

$p = $ppt.ActivePresentation

$updates = @(
    @{ Slide = 7;  Desc = "Introducing Greg Heffley and his family"; Label = "Character Introduction" },
    @{ Slide = 8;  Desc = "Greg's experiences navigating middle school"; Label = "School Life" },
    @{ Slide = 9;  Desc = "The ups and downs of Greg's relationship with his parents"; Label = "Family Dynamics" },
    @{ Slide = 10; Desc = "The complexities of friendships and peer pressure"; Label = "Friendship Challenges" },
    @{ Slide = 11; Desc = "Humorous stories that highlight the challenges of youth"; Label = "Anecdotes of Humor" },
    @{ Slide = 12; Desc = "Growth and lessons learned throughout the book"; Label = "Personal Growth" },
    @{ Slide = 13; Desc = "Reflecting on the journey and its implications for young readers."; Label = "Final Takeaways" }
)

foreach ($u in $updates) {
    $s = $p.Slides.Item($u.Slide)
    $s.Shapes.Item(8).TextFrame.TextRange.Text = $u.Desc
    $s.Shapes.Item(9).TextFrame.TextRange.Text = $u.Label
}
